$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.207.99"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "'3.458.25"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'578.83"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").Value = "'149.05"
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.479"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("D11").Value = "'0.409"
$ws.Range("E11").Value = "  +2.29%  "
$ws.Range("D12").Value = "'4.047.87"
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("E13").Value = "  +2.48%  "
$ws.Range("D14").Value = "'28.59"
$ws.Range("E14").Value = "  -4.39%  "
$ws.Range("D15").Value = "'3.454.42"
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").Value = "'63.161.51"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "'6.42"
$ws.Range("E18").Value = "  +2.30%  "
$ws.Range("D19").Value = "'14.45"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").Value = "'9.16"
$ws.Range("E20").Value = "  -3.41%  "
$ws.Range("D21").Value = "'387.04"
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("D22").Value = "'0.562"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "'74.47"
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'3.589.05"
$ws.Range("E25").Value = "  -1.56%  "
$ws.Range("E26").Value = "  -3.92%  "
$ws.Range("D27").Value = "'0.182"
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "'8.11"
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("E31").Value = "  -2.62%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "'23.38"
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("E34").Value = "  -6.16%  "
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("E36").Value = "  +2.80%  "
$ws.Range("D37").Value = "'7.06"
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("D38").Value = "'31.72"
$ws.Range("E38").Value = "  -2.85%  "
$ws.Range("D39").Value = "'170.08"
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("D40").Value = "'3.491.70"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("D41").Value = "'0.0770"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Value = "'0.792"
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("D43").Value = "'42.51"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("D45").Value = "'4.38"
$ws.Range("E45").Value = "  -2.95%  "
$ws.Range("E46").Value = "  -1.28%  "
$ws.Range("D47").Value = "'2.579.77"
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("D48").Value = "'2.31"
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("E49").Value = "  +1.51%  "
$ws.Range("D50").Value = "'22.63"
$ws.Range("E50").Value = "  -4.95%  "
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  -0.12%  "
